# Update column F (dSF) values to match repulled/recalculated data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = 4
$ws.Range("F5").Value = -3
$ws.Range("F6").Value = -5
$ws.Range("F7").Value = -4
$ws.Range("F8").Value = -6
$ws.Range("F9").Value = -5
$ws.Range("F11").Value = -5
$ws.Range("F12").Value = -4
$ws.Range("F19").Value = -14
$ws.Range("F20").Value = -7
$ws.Range("F22").Value = -9
